$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (B2 date changed, D2/E2 swapped, F2 value text changed)
$ws.Range("B2").Value = "2024-03-03 04:28:38"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = "40.000,00"

# New data rows 3-9
$data = @(
    @(2, "2024-03-18 10:28:46", "Pupuk", 100, 20, "200.000,00"),
    @(3, "2024-03-18 12:06:46", "Pupuk", 200, 200, "2.000.000,00"),
    @(4, "2024-03-18 12:05:31", "Pupuk", 1000, 0, "0,00"),
    @(5, "2024-03-18 12:29:28", "Pupuk", 1000, 0, "0,00"),
    @(6, "2024-03-18 12:30:25", "Pupuk", 1000, 0, "0,00"),
    @(7, "2024-03-18 12:33:31", "Pupuk", 800, 200, "2.000.000,00"),
    @(8, "2024-03-18 12:33:22", "Pupuk", 1000, 0, "0,00")
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}
